# Apply the risk-management tracker updates: mark several risks as
# Closed with updated execution summaries, bump the Testing row's
# target completion date, and turn on AutoFilter over the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Scope creep: execution summary update (stays Open)
$ws.Range("K7").Value = "Similar features have been condensed and combined. Feature prioritization and necessity is reevaluated at each sprint"

# Row 8 - Improper planning: execution summary + status -> Closed
$ws.Range("K8").Value = "Sprints have been planned realistically and work has been assigned fairly"
$ws.Range("L8").Value = "Closed"

# Rows 9-11 - Technology competence (React / unit testing / Git):
# execution summary text updated, status -> Closed
$ws.Range("K9").Value = "All team members have upskilled in their own time"
$ws.Range("L9").Value = "Closed"

$ws.Range("K10").Value = "All team members have upskilled in their own time"
$ws.Range("L10").Value = "Closed"

$ws.Range("K11").Value = "All team members have upskilled in their own time"
$ws.Range("L11").Value = "Closed"

# Row 13 - Not enough testing: target completion date pushed back
$ws.Range("I13").Value = "10/15/2024"
$ws.Range("I13").NumberFormat = "mm/dd/yyyy"
$ws.Range("I13").HorizontalAlignment = -4152

# Row 14 - Integration and deployment: execution summary + status -> Closed
$ws.Range("K14").Value = "dev branch is auto-deployed every time a branch is merged into it"
$ws.Range("L14").Value = "Closed"

# Row 15 - Hosting service: status -> Closed
$ws.Range("L15").Value = "Closed"

# Turn on the AutoFilter over the data table (header row 2 through last
# data row 15, columns A:L)
$ws.Range("A2:L15").AutoFilter()
